$wb = $excel.ActiveWorkbook

# --- Sheet "Notes" ---
$notes = $wb.Worksheets.Item("Notes")

# Drop the now-obsolete "Or remove U49-2..." line entirely (row 5) - everything below
# shifts up by one row
$notes.Rows.Item(5).Delete()

# Append the new CAN-communication jumper instructions as action item 3
$notes.Range("A8").Value = 3
$notes.Range("B8").Value = "CAN Communication jumpers"
$notes.Range("C9").Value = "Cut trace on W6, W7, W9 (solder side)"
$notes.Range("C10").Value = "Install jumpers in W4, W8, W10"

# Update wording of the two fault-fiber notes (RF / U48-6 references removed ahead of
# Maarten's fault-handling rewrite)
$notes.Range("B3").Value = 'Since there are no status fibers from PFN or Gun, the PIC will need to send status to clear these "faults"'
$notes.Range("C4").Value = "Either need to lift the outputs of the schmidt triggers (U48-2, U48-8) "

$notes.Range("C12").Select()

# --- Make "Notes" the active/selected tab (previously "Timer") ---
$notes.Activate()
